$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.404.29'
$ws.Range("E2").Value = '  +2.26%  '
$ws.Range("D3").Value = '1.461.49'
$ws.Range("E3").Value = '  +3.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.89%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9505'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '274.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3652'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3061'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.80'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06576'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.96'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.136'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001023'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").Value = '1.461.35'
$ws.Range("E17").Value = '  +3.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9701'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05818'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.428'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.24%  '
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = '20.412.91'
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.070'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.77%  '
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("D29").Value = '1.615.03'
$ws.Range("E29").Value = '  +2.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.852'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.884'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.07879'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7893'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.516'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05710'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.142'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.672'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.07%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02025'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9581'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.20%  '
$ws.Range("E41").Value = '  -1.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.462'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1856'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5252'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.487'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '117.13'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5139'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.745'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06418'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9927'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.73%  '
